$wb = $excel.ActiveWorkbook

# --- Update the "Conversión del día" text block on sheet "Hoja1" (cell A1) ---
$ws1 = $wb.Worksheets.Item("Hoja1")
$cell = $ws1.Range("A1")
$text = $cell.Value2
$text = $text -replace [regex]::Escape("1000 Bs = 2.68 = 9929.96 pesos"), "1000 Bs = 2.65 = 9803.77 pesos"
$text = $text -replace [regex]::Escape("9929.96 pesos = 2.67 = 944.92 Bs"), "9803.77 pesos = 2.64 = 929.7 Bs"
$cell.Value2 = $text

# --- Update the updated rates on sheet "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")
$ws2.Range("N10").Value = 376.998
$ws2.Range("O10").Value = 3696
$ws2.Range("O12").Value = 352.77
